# First run on dev env: switch the Parameters sheet's ENV to DEV and bump
# threadCount from 2 to 4.

$wb = $excel.ActiveWorkbook
$scenarios = $wb.Worksheets.Item("Scenarios")
$ws = $wb.Worksheets.Item("Parameters")

# ENV (row 8) -> DEV (baseUrl formula in B7 will recalc off this)
$ws.Range("B8").Value = "DEV"

# threadCount (row 13) -> 4
$ws.Range("B13").Value = 4

# Leave the selection on the cell the user last touched
$ws.Range("B13").Select()

# Restore the originally active sheet/tab (Scenarios) so only the
# Parameters sheet's own cursor position moves.
$scenarios.Activate()
